# Add a new "portNo" data point to the DeviceCapabilities sheet
# (mirrors the other header/value columns already present there).
$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("DeviceCapabilities")
$ws2.Range("H1").Value = "portNo"
$ws2.Range("H1").Font.Bold = $true
$ws2.Range("H2").Value = 4723

# Leave selection on A2 (mirrors how the sheet was left after the edit).
$ws2.Activate()
$ws2.Range("A2").Select()

$ws1 = $wb.Worksheets.Item("eBayTestData")
$ws1.Activate()
$ws1.Range("A2").Select()
